$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.129.43'

$ws.Range('D3').Value = '1.564.23'
$ws.Range('E3').Value = '  -1.98%  '

$ws.Range('E4').Value = '  -0.26%  '

$ws.Range('D5').Value = '''206.22'
$ws.Range('E5').Value = '  -1.41%  '

$ws.Range('E6').Value = '  -3.19%  '

$ws.Range('E7').Value = '  -0.23%  '

$ws.Range('D8').Value = '''22.24'
$ws.Range('E8').Value = '  -0.60%  '

$ws.Range('E9').Value = '  -2.55%  '

$ws.Range('E10').Value = '  -0.64%  '

$ws.Range('D11').Value = '''0.0861'
$ws.Range('E11').Value = '  -0.83%  '

$ws.Range('D12').Value = '1.784.99'
$ws.Range('E12').Value = '  -2.08%  '

$ws.Range('D13').Value = '1.565.36'
$ws.Range('E13').Value = '  -1.88%  '

$ws.Range('E14').Value = '  -2.39%  '

$ws.Range('E15').Value = '  -3.14%  '

$ws.Range('D16').Value = '''62.95'
$ws.Range('E16').Value = '  -0.86%  '

$ws.Range('D17').Value = '27.125.37'
$ws.Range('E17').Value = '  -2.30%  '

$ws.Range('D18').Value = '''215.11'
$ws.Range('E18').Value = '  -2.10%  '

$ws.Range('E19').Value = '  -1.96%  '

$ws.Range('E20').Value = '  -1.91%  '

$ws.Range('E21').Value = '  -0.15%  '

$ws.Range('E22').Value = '  -0.95%  '

$ws.Range('E23').Value = '  -4.75%  '

$ws.Range('E24').Value = '  +0.31%  '

$ws.Range('D25').Value = '''151.59'
$ws.Range('E25').Value = '  -1.60%  '

$ws.Range('E26').Value = '  -7.89%  '

$ws.Range('D27').Value = '''14.93'
$ws.Range('E27').Value = '  -1.68%  '

$ws.Range('E28').Value = '  -0.20%  '

$ws.Range('E29').Value = '  -1.72%  '

$ws.Range('E30').Value = '  -1.96%  '

$ws.Range('D31').Value = '''0.0461'
$ws.Range('E31').Value = '  -2.56%  '

$ws.Range('D32').Value = '''3.17'
$ws.Range('E32').Value = '  -2.19%  '

$ws.Range('D33').Value = '1.388.80'
$ws.Range('E33').Value = '  +0.87%  '

$ws.Range('D34').Value = '''2.92'
$ws.Range('E34').Value = '  -1.86%  '

$ws.Range('E35').Value = '  +0.04%  '

$ws.Range('D36').Value = '''2.28'
$ws.Range('E36').Value = '  -2.14%  '

$ws.Range('D37').Value = '''0.940'
$ws.Range('E37').Value = '  -3.47%  '

$ws.Range('E38').Value = '  -1.99%  '

$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '''0.812'
$ws.Range('E39').Value = '  -2.02%  '

$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').Value = '''0.517'
$ws.Range('E40').Value = '  -3.87%  '

$ws.Range('E41').Value = '  -0.15%  '

$ws.Range('D42').Value = '''0.989'
$ws.Range('E42').Value = '  +1.83%  '

$ws.Range('D43').Value = '''1.79'
$ws.Range('E43').Value = '  +1.92%  '

$ws.Range('D44').Value = '''63.50'
$ws.Range('E44').Value = '  -1.72%  '

$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '''5.28'
$ws.Range('E45').Value = '  +0.45%  '

$ws.Range('B46').Value = 'MXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D46').Value = '''2.16'
$ws.Range('E46').Value = '  -0.76%  '

$ws.Range('D47').Value = '1.698.42'
$ws.Range('E47').Value = '  -2.07%  '

$ws.Range('D48').Value = '''85.62'
$ws.Range('E48').Value = '  -1.21%  '

$ws.Range('D49').Value = '0.0₇0982'
$ws.Range('E49').Value = '  -3.16%  '

$ws.Range('E50').Value = '  -1.07%  '

$ws.Range('D51').Value = '''0.0945'
$ws.Range('E51').Value = '  -2.34%  '
